# 自动更新Excel文件 - 2025-12-17 23:13:33
# This script advances the "days remaining" tracker in column E by one day
# (the daily update run), rolling the "start date" (column F) forward to a
# new cycle whenever the remaining-day counter would drop to zero or below.
#
#   E (剩余/remaining) = D (总天/total days) - daysElapsed(F -> today)
#
# Before this run "today" was 2025-12-15 + ... i.e. derived from the data as
# 2025-12-17; this run moves "today" forward to 2025-12-18. Whenever the
# newly computed remaining value would be <= 0, the row is treated as having
# completed its cycle: it restarts today (F = new "today") with a fresh full
# remaining count (E = D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldToday = [datetime]::ParseExact("20251217", "yyyyMMdd", [System.Globalization.CultureInfo]::InvariantCulture)
$newToday = [datetime]::ParseExact("20251218", "yyyyMMdd", [System.Globalization.CultureInfo]::InvariantCulture)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fText = [string][int]$fVal

    $fDate = $null
    try {
        $fDate = [datetime]::ParseExact($fText, "yyyyMMdd", [System.Globalization.CultureInfo]::InvariantCulture)
    } catch {
        # Unparseable start date (corrupted data) - leave the row untouched.
        continue
    }

    $elapsed = [int]($oldToday.ToOADate() - $fDate.ToOADate())
    $newRemaining = [int]$dVal - $elapsed - 1

    if ($newRemaining -le 0) {
        # Cycle complete: restart today with a full remaining count.
        $eCell.Value = [int]$dVal
        $fCell.Value = [int]$newToday.ToString("yyyyMMdd")
    } else {
        $eCell.Value = $newRemaining
    }
}
